$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 6
$ws.Range("H2").Value = 6

# Row 4
$ws.Range("E4").Value = 14

# Row 10
$ws.Range("E10").Value = 12

# Row 11
$ws.Range("E11").Value = 9

# Row 12
$ws.Range("E12").Value = 16

# Row 15
$ws.Range("E15").Value = 58

# Row 16
$ws.Range("E16").Value = 221
$ws.Range("F16").Value = 56
$ws.Range("H16").Value = 56

# Row 18
$ws.Range("E18").Value = 57
